$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.104.91"
$ws.Range("E2").Value = "  +5.55%  "
$ws.Range("D3").Value = "2.263.93"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'302.46"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").Value = "'93.15"
$ws.Range("E6").Value = "  +7.39%  "
$ws.Range("E7").Value = "  +4.33%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +4.24%  "
$ws.Range("D10").Value = "'32.88"
$ws.Range("E10").Value = "  +8.55%  "
$ws.Range("D11").Value = "'54.88"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "2.614.60"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "'14.17"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "2.269.80"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("E18").Value = "  +4.06%  "
$ws.Range("D19").Value = "42.007.36"
$ws.Range("E19").Value = "  +5.64%  "
$ws.Range("D20").Value = "'12.22"
$ws.Range("E20").Value = "  +9.18%  "
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("E22").Value = "  +3.86%  "
$ws.Range("D23").Value = "'67.41"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'242.46"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").Value = "'23.98"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("E30").Value = "  +5.64%  "
$ws.Range("D31").Value = "'34.34"
$ws.Range("E31").Value = "  +7.88%  "
$ws.Range("D32").Value = "'158.68"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +5.02%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("D36").Value = "'3.09"
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("D39").Value = "'16.65"
$ws.Range("E39").Value = "  +9.59%  "
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").Value = "'20.29"
$ws.Range("E43").Value = "  +14.86%  "
$ws.Range("D44").Value = "2.057.22"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  +4.91%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "'2.93"
$ws.Range("E47").Value = "  +8.37%  "
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("D49").Value = "2.486.83"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.14"
$ws.Range("E51").Value = "  +4.97%  "
